$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name (B) updates ---
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('B10').Value = 'Polygon'
$ws.Range('B11').Value = 'Solana'
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('B15').Value = 'TRON'
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('B19').Value = 'Dai'
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('B26').Value = 'Monero'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('B31').Value = 'Stellar'
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('B36').Value = 'Frax'
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('B39').Value = 'VeChain'
$ws.Range('B40').Value = 'MXToken'
$ws.Range('B41').Value = 'Algorand'
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('B44').Value = 'Aptos'
$ws.Range('B45').Value = 'Cronos'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('B47').Value = 'Quant'
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('B51').Value = 'Aave'

# --- Link (C) updates ---
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

# --- Price (D) updates: values that are NOT valid numbers stay text automatically ---
$ws.Range('D2').Value = '27.688.30'
$ws.Range('D3').Value = '1.844.64'
$ws.Range('D12').Value = '1.858.24'
$ws.Range('D21').Value = '27.903.86'
$ws.Range('D24').Value = '2.099.74'

# --- Price (D) updates: values that look numeric must be forced to Text format ---
# so they remain stored as text (matching the source data type) instead of being
# auto-converted to numbers by Excel.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.24'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4293'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3637'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07325'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8765'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.71'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.344'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.517'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06952'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.005'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.51'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008998'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.38'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.978'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.30'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.992'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.75'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.57'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '119.65'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.211'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.880'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08889'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7533'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.529'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.960'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.118'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.002'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.108'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05428'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01935'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.827'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1668'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5070'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.620'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.371'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06543'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.39'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.76'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4650'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.637'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.69'

# --- Volume(1h) (E) updates ---
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('E12').Value = '  +2.75%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +3.36%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('E35').Value = '  +1.27%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('E43').Value = '  -2.94%  '
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('E51').Value = '  +0.48%  '
